$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.261.11'
$ws.Range('E2').Value = '  +0.69%  '
$ws.Range('D3').Value = '2.025.24'
$ws.Range('E3').Value = '  +0.05%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'228.79"
$ws.Range('E5').Value = '  +1.15%  '
$ws.Range('E6').Value = '  +0.24%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = "'56.20"
$ws.Range('E8').Value = '  +2.16%  '
$ws.Range('E9').Value = '  -0.82%  '
$ws.Range('E10').Value = '  -1.09%  '
$ws.Range('D11').Value = "'0.103"
$ws.Range('E11').Value = '  -1.96%  '
$ws.Range('D12').Value = '2.324.89'
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('D13').Value = "'14.32"
$ws.Range('E13').Value = '  +0.26%  '
$ws.Range('D14').Value = "'20.22"
$ws.Range('E14').Value = '  -1.80%  '
$ws.Range('D15').Value = "'0.741"
$ws.Range('E15').Value = '  -0.56%  '
$ws.Range('E16').Value = '  +0.81%  '
$ws.Range('D17').Value = '2.030.59'
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('D18').Value = '37.248.55'
$ws.Range('E18').Value = '  +0.85%  '
$ws.Range('D19').Value = "'6.16"
$ws.Range('E19').Value = '  +2.19%  '
$ws.Range('D20').Value = "'69.02"
$ws.Range('E20').Value = '  +0.22%  '
$ws.Range('E21').Value = '  -1.42%  '
$ws.Range('D22').Value = "'222.91"
$ws.Range('E22').Value = '  -1.25%  '
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('E24').Value = '  +1.85%  '
$ws.Range('D25').Value = "'2.23"
$ws.Range('E25').Value = '  -1.50%  '
$ws.Range('D26').Value = "'163.71"
$ws.Range('E26').Value = '  -2.37%  '
$ws.Range('E27').Value = '  -2.98%  '
$ws.Range('D28').Value = "'0.129"
$ws.Range('E28').Value = '  +2.48%  '
$ws.Range('D29').Value = "'18.68"
$ws.Range('E29').Value = '  -0.59%  '
$ws.Range('D30').Value = "'1.31"
$ws.Range('E30').Value = '  -2.09%  '
$ws.Range('D31').Value = "'0.117"
$ws.Range('E31').Value = '  +0.40%  '
$ws.Range('D32').Value = "'4.45"
$ws.Range('E32').Value = '  -0.60%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = "'0.0605"
$ws.Range('E33').Value = '  -0.84%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').Value = "'2.01"
$ws.Range('E34').Value = '  +10.17%  '
$ws.Range('D35').Value = "'4.45"
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('E36').Value = '  -1.45%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').Value = "'3.21"
$ws.Range('E37').Value = '  +0.47%  '
$ws.Range('B38').Value = 'BinanceUSD'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D38').Value = "'0.999"
$ws.Range('E38').Value = '  -0.23%  '
$ws.Range('D39').Value = "'5.55"
$ws.Range('E39').Value = '  +2.02%  '
$ws.Range('B40').Value = 'FTXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D40').Value = "'4.49"
$ws.Range('E40').Value = '  +23.53%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '1.469.64'
$ws.Range('E41').Value = '  -1.99%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = "'0.0213"
$ws.Range('E42').Value = '  -2.51%  '
$ws.Range('D43').Value = "'2.81"
$ws.Range('E43').Value = '  -0.45%  '
$ws.Range('D44').Value = "'94.01"
$ws.Range('E44').Value = '  -1.36%  '
$ws.Range('E45').Value = '  -1.21%  '
$ws.Range('D46').Value = "'16.24"
$ws.Range('E46').Value = '  -3.75%  '
$ws.Range('E47').Value = '  -1.97%  '
$ws.Range('E48').Value = '  +0.64%  '
$ws.Range('D49').Value = "'7.15"
$ws.Range('E49').Value = '  -1.00%  '
$ws.Range('E50').Value = '  +1.03%  '
$ws.Range('D51').Value = '2.212.79'
$ws.Range('E51').Value = '  +0.03%  '
